# Update the "built on" timestamp embedded in the version string from
# "January 30 2026 16.19.47 EST" to "February 02 2026 12.49.33 EST"
# across the "About" sheet and the "Boundaries and methane sources" sheet.

$wb = $excel.ActiveWorkbook

$oldStamp = "January 30 2026 16.19.47 EST"
$newStamp = "February 02 2026 12.49.33 EST"

$wsAbout = $wb.Worksheets.Item("About")
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")

# --- "About" sheet ---
$valA2 = $wsAbout.Range("A2").Value()
$wsAbout.Range("A2").Value = $valA2.Replace($oldStamp, $newStamp)

$valA6 = $wsAbout.Range("A6").Value()
$wsAbout.Range("A6").Value = $valA6.Replace($oldStamp, $newStamp)

# --- "Boundaries and methane sources" sheet: column S, rows 2-8 ---
for ($row = 2; $row -le 8; $row++) {
    $cell = $wsData.Cells.Item($row, 19)  # column S = 19
    $val = $cell.Value()
    $cell.Value = $val.Replace($oldStamp, $newStamp)
}
